$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Rows 1-4 (header + the 白頭翁/麻雀 and 麻雀/雀榕 relations) are unchanged.
# Rows 5-9 get updated values, and rows 10-13 are newly appended, adding two
# new taxa - 五色鳥 (Taiwan barbet) and 臺灣海棗 (Taiwan date palm) - and
# their feeding relationships to the existing predator/prey list.
#
# Row 8 is written before row 5 so that the brand-new shared string "五色鳥"
# gets registered ahead of "臺灣海棗" (matching the order newly-added taxa
# were appended to the workbook's shared-string table).

$ws.Cells.Item(8, 1).Value = "五色鳥"
$ws.Cells.Item(8, 2).Value = "臺灣海棗"
$ws.Cells.Item(8, 3).Value = "吃"

$ws.Cells.Item(5, 1).Value = "麻雀"
$ws.Cells.Item(5, 2).Value = "臺灣海棗"
$ws.Cells.Item(5, 3).Value = "吃"

$ws.Cells.Item(6, 1).Value = "雀榕"
$ws.Cells.Item(6, 2).Value = "麻雀"
$ws.Cells.Item(6, 3).Value = "被吃"

$ws.Cells.Item(7, 1).Value = "紋翼畫眉"
$ws.Cells.Item(7, 2).Value = "姑婆芋"
$ws.Cells.Item(7, 3).Value = "吃"

$ws.Cells.Item(9, 1).Value = "綠背斜紋天蛾"
$ws.Cells.Item(9, 2).Value = "密毛魔芋"
$ws.Cells.Item(9, 3).Value = "吃"

$ws.Cells.Item(10, 1).Value = "姑婆芋"
$ws.Cells.Item(10, 2).Value = "紋翼畫眉"
$ws.Cells.Item(10, 3).Value = "被吃"

$ws.Cells.Item(11, 1).Value = "密毛魔芋"
$ws.Cells.Item(11, 2).Value = "綠背斜紋天蛾"
$ws.Cells.Item(11, 3).Value = "被吃"

$ws.Cells.Item(12, 1).Value = "臺灣海棗"
$ws.Cells.Item(12, 2).Value = "麻雀"
$ws.Cells.Item(12, 3).Value = "被吃"

$ws.Cells.Item(13, 1).Value = "臺灣海棗"
$ws.Cells.Item(13, 2).Value = "五色鳥"
$ws.Cells.Item(13, 3).Value = "被吃"

# The newly-created rows (10-13) don't automatically inherit the left-aligned
# body-row formatting used by the rest of the table, so re-apply it
# explicitly (matches the existing style used by A2:C9).
$ws.Range("A10:C13").HorizontalAlignment = -4131
